$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2279
$ws.Range("F3").Value = 352
$ws.Range("F4").Value = 184
$ws.Range("F5").Value = 187
$ws.Range("F6").Value = 348
$ws.Range("F10").Value = 678
$ws.Range("F11").Value = 372
$ws.Range("F12").Value = 64
$ws.Range("F15").Value = 6394
$ws.Range("F16").Value = 192
$ws.Range("F17").Value = 17
$ws.Range("F18").Value = 36
$ws.Range("F20").Value = 144
$ws.Range("F21").Value = 113
$ws.Range("F23").Value = 101
$ws.Range("F26").Value = 109

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 184
$ws.Range("F7").Value = 210
$ws.Range("F8").Value = 2836
$ws.Range("F10").Value = 25
$ws.Range("F16").Value = 2552

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 48
$ws.Range("F4").Value = 412
$ws.Range("F5").Value = 175

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 48
$ws.Range("F6").Value = 2279
$ws.Range("F7").Value = 412
$ws.Range("F8").Value = 352
$ws.Range("F9").Value = 184
$ws.Range("F10").Value = 187
$ws.Range("F11").Value = 348
$ws.Range("F15").Value = 184
$ws.Range("F16").Value = 175
$ws.Range("F19").Value = 678
$ws.Range("F20").Value = 372
$ws.Range("F21").Value = 64
$ws.Range("F24").Value = 6398
$ws.Range("F25").Value = 210
$ws.Range("F26").Value = 2836
$ws.Range("F28").Value = 25
$ws.Range("F30").Value = 192
$ws.Range("F31").Value = 17
$ws.Range("F32").Value = 36
$ws.Range("F36").Value = 144
$ws.Range("F37").Value = 113
$ws.Range("F41").Value = 101
$ws.Range("F44").Value = 109
$ws.Range("F45").Value = 2552
